$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.656.54"
$ws.Range("E2").Value = "  +3.88%  "
$ws.Range("D3").Value = "1.919.74"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.19"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.42"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.92"
$ws.Range("E9").Value = "  +10.40%  "
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0766"
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0999"
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.61"
$ws.Range("E13").Value = "  +8.68%  "
$ws.Range("E14").Value = "  +3.90%  "
$ws.Range("D15").Value = "2.197.37"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("E16").Value = "  +4.74%  "
$ws.Range("D17").Value = "1.922.42"
$ws.Range("E17").Value = "  +2.46%  "
$ws.Range("D18").Value = "36.674.63"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.28"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +5.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "251.75"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.27"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  +5.25%  "
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.00"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.83"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.78"
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.55"
$ws.Range("E31").Value = "  +7.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0617"
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.37"
$ws.Range("E34").Value = "  +5.98%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0850"
$ws.Range("E36").Value = "  +14.65%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  -12.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.877"
$ws.Range("E38").Value = "  +4.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.83"
$ws.Range("E39").Value = "  +47.06%  "
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.72"
$ws.Range("E41").Value = "  +11.01%  "
$ws.Range("E42").Value = "  +6.05%  "
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").Value = "1.339.61"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("E46").Value = "  +6.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0814"
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("E50").Value = "  +3.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.05"
$ws.Range("E51").Value = "  +2.41%  "
